$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Cell value edits (also rewrites the shared-strings table) ---
# Order matters: new shared strings are appended in the order the cells are
# written, and the target file expects "pistolpete" (88) before ":flag_kr:" (89).
$ws.Cells.Item(7, 2).Value = "pistolpete"
$ws.Cells.Item(4, 6).Value = ":flag_kr:"

# --- Column visibility / width ---
$ws.Columns.Item(3).Hidden = $true
$ws.Columns.Item(4).Hidden = $true
$ws.Columns.Item(5).Hidden = $true
$ws.Columns.Item(7).Hidden = $true
$ws.Columns.Item(8).Hidden = $true
$ws.Columns.Item(9).ColumnWidth = 5.45  # nearest representable value to the raw 6.28515625 OOXML width
$ws.Columns.Item(10).Hidden = $true
$ws.Columns.Item(11).Hidden = $true
$ws.Range($ws.Columns.Item(12), $ws.Columns.Item(22)).Hidden = $true

# --- Sheet view: zoom + selection ---
$win = $excel.ActiveWindow
$win.Zoom = 115
$ws.Range("F4").Select()
